$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "General"
$ws.Range("D2").Value = "'316.62"
$ws.Range("E2").NumberFormat = "General"
$ws.Range("E2").Value = "'1.44%"
$ws.Range("D3").NumberFormat = "General"
$ws.Range("D3").Value = "'37.84"
$ws.Range("E3").NumberFormat = "General"
$ws.Range("E3").Value = "'1.03%"
$ws.Range("D4").NumberFormat = "General"
$ws.Range("D4").Value = "'5.187"
$ws.Range("E4").NumberFormat = "General"
$ws.Range("E4").Value = "'1.27%"
$ws.Range("D5").NumberFormat = "General"
$ws.Range("D5").Value = "'0.07989"
$ws.Range("E5").NumberFormat = "General"
$ws.Range("E5").Value = "'1.64%"
$ws.Range("B6").Value = "GateToken"
$ws.Range("C6").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D6").NumberFormat = "General"
$ws.Range("D6").Value = "'4.490"
$ws.Range("E6").NumberFormat = "General"
$ws.Range("E6").Value = "'1.55%"
$ws.Range("B7").Value = "KuCoinToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
$ws.Range("D7").NumberFormat = "General"
$ws.Range("D7").Value = "'8.541"
$ws.Range("E7").NumberFormat = "General"
$ws.Range("E7").Value = "'3.44%"
$ws.Range("B8").Value = "FTXToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D8").NumberFormat = "General"
$ws.Range("D8").Value = "'1.928"
$ws.Range("E8").NumberFormat = "General"
$ws.Range("E8").Value = "'0.62%"
$ws.Range("B9").Value = "BTSEToken"
$ws.Range("C9").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D9").NumberFormat = "General"
$ws.Range("D9").Value = "'2.963"
$ws.Range("E9").NumberFormat = "General"
$ws.Range("E9").Value = "'0.37%"
$ws.Range("B10").Value = "MXToken"
$ws.Range("C10").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D10").NumberFormat = "General"
$ws.Range("D10").Value = "'0.9445"
$ws.Range("E10").NumberFormat = "General"
$ws.Range("E10").Value = "'2.75%"
$ws.Range("B11").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C11").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D11").NumberFormat = "General"
$ws.Range("D11").Value = "'0.1299"
$ws.Range("E11").NumberFormat = "General"
$ws.Range("E11").Value = "'8.43%"
$ws.Range("B12").Value = "WazirX"
$ws.Range("C12").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D12").NumberFormat = "General"
$ws.Range("D12").Value = "'0.1942"
$ws.Range("E12").NumberFormat = "General"
$ws.Range("E12").Value = "'1.58%"
$ws.Range("B13").Value = "MandalaExchangeToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D13").NumberFormat = "General"
$ws.Range("D13").Value = "'0.09105"
$ws.Range("E13").NumberFormat = "General"
$ws.Range("E13").Value = "'0.98%"
$ws.Range("B14").Value = "BitrueCoin"
$ws.Range("C14").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D14").NumberFormat = "General"
$ws.Range("D14").Value = "'0.03401"
$ws.Range("E14").NumberFormat = "General"
$ws.Range("E14").Value = "'1.56%"
$ws.Range("B15").Value = "BitMartToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D15").NumberFormat = "General"
$ws.Range("D15").Value = "'0.09528"
$ws.Range("E15").NumberFormat = "General"
$ws.Range("E15").Value = "'-0.65%"
$ws.Range("B16").Value = "BitForexToken"
$ws.Range("C16").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D16").NumberFormat = "General"
$ws.Range("D16").Value = "'0.001396"
$ws.Range("E16").NumberFormat = "General"
$ws.Range("E16").Value = "'1.18%"
$ws.Range("B17").Value = "TigerCash"
$ws.Range("C17").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D17").NumberFormat = "General"
$ws.Range("D17").Value = "'0.005950"
$ws.Range("E17").NumberFormat = "General"
$ws.Range("E17").Value = "'4.01%"
$ws.Range("B18").Value = "LEO"
$ws.Range("C18").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D18").NumberFormat = "General"
$ws.Range("D18").Value = "'3.432"
$ws.Range("E18").NumberFormat = "General"
$ws.Range("E18").Value = "'-2.78%"
$ws.Range("D19").NumberFormat = "General"
$ws.Range("D19").Value = "'0.3516"
$ws.Range("E19").NumberFormat = "General"
$ws.Range("E19").Value = "'2.12%"
$ws.Range("D20").NumberFormat = "General"
$ws.Range("D20").Value = "'6.601"
$ws.Range("E20").NumberFormat = "General"
$ws.Range("E20").Value = "'26.60%"
$ws.Range("D21").NumberFormat = "General"
$ws.Range("D21").Value = "'0.1299"
$ws.Range("E21").NumberFormat = "General"
$ws.Range("E21").Value = "'1.20%"
$ws.Range("D22").NumberFormat = "General"
$ws.Range("D22").Value = "'0.2424"
$ws.Range("E22").NumberFormat = "General"
$ws.Range("E22").Value = "'-6.37%"
$ws.Range("D23").NumberFormat = "General"
$ws.Range("D23").Value = "'0.04378"
$ws.Range("E23").NumberFormat = "General"
$ws.Range("E23").Value = "'0.60%"
$ws.Range("E24").NumberFormat = "General"
$ws.Range("E24").Value = "'-1.42%"
$ws.Range("D25").NumberFormat = "General"
$ws.Range("D25").Value = "'0.004261"
$ws.Range("E25").NumberFormat = "General"
$ws.Range("E25").Value = "'-8.52%"
$ws.Range("D26").NumberFormat = "General"
$ws.Range("D26").Value = "'0.0001326"
$ws.Range("E26").NumberFormat = "General"
$ws.Range("E26").Value = "'-1.97%"
$ws.Range("D27").NumberFormat = "General"
$ws.Range("D27").Value = "'0.0003986"
$ws.Range("E27").NumberFormat = "General"
$ws.Range("E27").Value = "'-0.08%"
$ws.Range("D39").NumberFormat = "General"
$ws.Range("D39").Value = "'0.02379"
$ws.Range("E39").NumberFormat = "General"
$ws.Range("E39").Value = "'5.37%"
$ws.Range("D40").NumberFormat = "General"
$ws.Range("D40").Value = "'0.05158"
$ws.Range("E40").NumberFormat = "General"
$ws.Range("E40").Value = "'1.93%"
$ws.Range("D41").NumberFormat = "General"
$ws.Range("D41").Value = "'0.007611"
$ws.Range("E41").NumberFormat = "General"
$ws.Range("E41").Value = "'1.97%"
$ws.Range("E42").NumberFormat = "General"
$ws.Range("E42").Value = "'3.82%"
$ws.Range("D43").NumberFormat = "General"
$ws.Range("D43").Value = "'0.008562"
$ws.Range("E43").NumberFormat = "General"
$ws.Range("E43").Value = "'-5.39%"
$ws.Range("D44").NumberFormat = "General"
$ws.Range("D44").Value = "'0.002104"
$ws.Range("E44").NumberFormat = "General"
$ws.Range("E44").Value = "'8.33%"
$ws.Range("D45").NumberFormat = "General"
$ws.Range("D45").Value = "'0.008748"
$ws.Range("E45").NumberFormat = "General"
$ws.Range("E45").Value = "'-5.49%"
$ws.Range("D46").NumberFormat = "General"
$ws.Range("D46").Value = "'0.00006492"
$ws.Range("E46").NumberFormat = "General"
$ws.Range("E46").Value = "'-1.02%"
$ws.Range("E47").NumberFormat = "General"
$ws.Range("E47").Value = "'-0.08%"
$ws.Range("D48").NumberFormat = "General"
$ws.Range("D48").Value = "'0.002864"
$ws.Range("E48").NumberFormat = "General"
$ws.Range("E48").Value = "'-14.77%"
$ws.Range("D49").NumberFormat = "General"
$ws.Range("D49").Value = "'0.001688"
$ws.Range("E49").NumberFormat = "General"
$ws.Range("E49").Value = "'68.81%"
$ws.Range("D50").NumberFormat = "General"
$ws.Range("D50").Value = "'0.00002098"
$ws.Range("E50").NumberFormat = "General"
$ws.Range("E50").Value = "'-0.08%"
$ws.Range("D51").NumberFormat = "General"
$ws.Range("D51").Value = "'0.0001998"
$ws.Range("E51").NumberFormat = "General"
$ws.Range("E51").Value = "'-0.08%"
